$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.422.69'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '3.388.03'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.21'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.38'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.474'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('E9').Value = '  +1.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.122'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.387'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('D12').Value = '3.973.02'
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.34'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('D15').Value = '3.399.86'
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000170'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('D17').Value = '61.523.98'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.15'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.65'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.97'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '388.42'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.18'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.555'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000112'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.39%  '
$ws.Range('E26').Value = '  +6.60%  '
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.29'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.03'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('E30').Value = '  -0.59%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.37'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.39'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.93'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.48%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '167.59'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.03'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('D37').Value = '3.428.42'
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.47'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0769'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '25.77'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -10.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.777'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.44'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.66'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.72%  '
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('D45').Value = '2.448.26'
$ws.Range('E45').Value = '  -2.18%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.71'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.71%  '
$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.70'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.97%  '
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0262'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.03'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.206'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.99%  '
